$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds "2021-01-09" which Excel would otherwise auto-parse as a
# date serial; force it to stay plain text like the rest of the sheet by
# pre-formatting as Text, then clearing the style back to Normal once the
# value is set (matches the un-styled cells produced by the original data).
$dateCols = @("A75","A76","A77","A78","A79","A80","A81","A82","A83","A84")
foreach ($addr in $dateCols) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A75").Value = "2021-01-09"
$ws.Range("B75").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("C75").Value = 3333.57
$ws.Range("D75").Value = 1033.16
$ws.Range("E75").Value = 2300.41
$ws.Range("F75").Value = "KNN"
$ws.Range("J75").Value = 2171.66
$ws.Range("K75").Value = 70.61

$ws.Range("A76").Value = "2021-01-09"
$ws.Range("B76").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D76").Value = 1049.43
$ws.Range("F76").Value = "KNN"

$ws.Range("A77").Value = "2021-01-09"
$ws.Range("B77").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D77").Value = 874.74
$ws.Range("F77").Value = "KNN"

$ws.Range("A78").Value = "2021-01-09"
$ws.Range("B78").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D78").Value = 771.86
$ws.Range("F78").Value = "KNN"

$ws.Range("A79").Value = "2021-01-09"
$ws.Range("B79").Value = "07 Feb -- 13 Feb 2021"
$ws.Range("D79").Value = 821.02
$ws.Range("F79").Value = "KNN"

$ws.Range("A80").Value = "2021-01-09"
$ws.Range("B80").Value = "14 Feb -- 20 Feb 2021"
$ws.Range("D80").Value = 816.24
$ws.Range("F80").Value = "KNN"

$ws.Range("A81").Value = "2021-01-09"
$ws.Range("B81").Value = "21 Feb -- 27 Feb 2021"
$ws.Range("D81").Value = 1062.92
$ws.Range("F81").Value = "KNN"

$ws.Range("A82").Value = "2021-01-09"
$ws.Range("B82").Value = "28 Feb -- 06 Mar 2021"
$ws.Range("D82").Value = 1051.94
$ws.Range("F82").Value = "KNN"

$ws.Range("A83").Value = "2021-01-09"
$ws.Range("B83").Value = "07 Mar -- 13 Mar 2021"
$ws.Range("D83").Value = 778.59
$ws.Range("F83").Value = "KNN"

$ws.Range("A84").Value = "2021-01-09"
$ws.Range("B84").Value = "14 Mar -- 20 Mar 2021"
$ws.Range("D84").Value = 859.84
$ws.Range("F84").Value = "KNN"

# Clear the Text number format back to Normal/General on column A so the
# new rows end up with no style index, matching the rest of the sheet.
foreach ($addr in $dateCols) {
  $ws.Range($addr).Style = "Normal"
}
